# Auto-generated: refresh market-price derived columns (H-N) across all job sheets
# per the scheduled-runner data update described in the commit.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H18").Value = 1265.25
$ws.Range("I18").Value = 1265.25
$ws.Range("K18").Value = 1265.25
$ws.Range("M18").Value = -981.25
$ws.Range("H70").Value = 1671.25
$ws.Range("I70").Value = 1166.6666
$ws.Range("J70").Value = 1974
$ws.Range("K70").Value = 3499.9998
$ws.Range("L70").Value = 5922
$ws.Range("M70").Value = -3229.9998
$ws.Range("N70").Value = -6462
$ws.Range("H73").Value = 1671.25
$ws.Range("I73").Value = 1166.6666
$ws.Range("J73").Value = 1974
$ws.Range("K73").Value = 3499.9998
$ws.Range("L73").Value = 5922
$ws.Range("M73").Value = -2563.9998
$ws.Range("N73").Value = -7794
$ws.Range("H98").Value = 1912.1072
$ws.Range("I98").Value = 1884.9584
$ws.Range("K98").Value = 1884.9584
$ws.Range("M98").Value = -386.9584
$ws.Range("H112").Value = 2766.6667
$ws.Range("I112").Value = 1850
$ws.Range("J112").Value = 2863.158
$ws.Range("K112").Value = 5550
$ws.Range("L112").Value = 8589.474
$ws.Range("M112").Value = -4442
$ws.Range("N112").Value = -10805.474
$ws.Range("H122").Value = 1912.1072
$ws.Range("I122").Value = 1884.9584
$ws.Range("K122").Value = 5654.8752
$ws.Range("M122").Value = -3204.8752
$ws.Range("H137").Value = 1085.4762
$ws.Range("I137").Value = 806.3333
$ws.Range("J137").Value = 1783.3334
$ws.Range("K137").Value = 2418.9999
$ws.Range("L137").Value = 5350.0002
$ws.Range("M137").Value = 131.0001000000002
$ws.Range("N137").Value = -10450.0002
$ws.Range("H138").Value = 3933.6594
$ws.Range("I138").Value = 2630.9412
$ws.Range("J138").Value = 4232.9326
$ws.Range("K138").Value = 7892.823600000001
$ws.Range("L138").Value = 12698.7978
$ws.Range("M138").Value = -2752.823600000001
$ws.Range("N138").Value = -22978.7978

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 17535.812
$ws.Range("I32").Value = 12842.359
$ws.Range("J32").Value = 69834.28999999999
$ws.Range("K32").Value = 12842.359
$ws.Range("L32").Value = 69834.28999999999
$ws.Range("M32").Value = -12555.359
$ws.Range("N32").Value = -70408.28999999999
$ws.Range("H61").Value = 1112.1333
$ws.Range("I61").Value = 1114.6666
$ws.Range("J61").Value = 1102
$ws.Range("K61").Value = 1114.6666
$ws.Range("L61").Value = 1102
$ws.Range("M61").Value = -902.6666
$ws.Range("N61").Value = -1526
$ws.Range("H74").Value = 1210.3611
$ws.Range("I74").Value = 1206.3448
$ws.Range("J74").Value = 1227
$ws.Range("K74").Value = 1206.3448
$ws.Range("L74").Value = 1227
$ws.Range("M74").Value = -332.3448000000001
$ws.Range("N74").Value = -2975
$ws.Range("H77").Value = 1210.3611
$ws.Range("I77").Value = 1206.3448
$ws.Range("J77").Value = 1227
$ws.Range("K77").Value = 6031.724
$ws.Range("L77").Value = 6135
$ws.Range("M77").Value = -1663.724
$ws.Range("N77").Value = -14871
$ws.Range("H136").Value = 1112.1333
$ws.Range("I136").Value = 1114.6666
$ws.Range("J136").Value = 1102
$ws.Range("K136").Value = 3343.9998
$ws.Range("L136").Value = 3306
$ws.Range("M136").Value = -793.9998000000001
$ws.Range("N136").Value = -8406

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 13910.263
$ws.Range("I134").Value = 1198.2206
$ws.Range("J134").Value = 85945.164
$ws.Range("K134").Value = 3594.6618
$ws.Range("L134").Value = 257835.492
$ws.Range("M134").Value = -1059.6618
$ws.Range("N134").Value = -262905.492
$ws.Range("H141").Value = 50210.1
$ws.Range("J141").Value = 50210.1
$ws.Range("L141").Value = 50210.1
$ws.Range("N141").Value = -60570.1

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 348.5
$ws.Range("I22").Value = 180.16667
$ws.Range("J22").Value = 601
$ws.Range("K22").Value = 180.16667
$ws.Range("L22").Value = 601
$ws.Range("M22").Value = 169.83333
$ws.Range("N22").Value = -1301
$ws.Range("H31").Value = 2247.8462
$ws.Range("I31").Value = 1855.5682
$ws.Range("J31").Value = 4405.375
$ws.Range("K31").Value = 1855.5682
$ws.Range("L31").Value = 4405.375
$ws.Range("M31").Value = -1560.5682
$ws.Range("N31").Value = -4995.375
$ws.Range("H34").Value = 2247.8462
$ws.Range("I34").Value = 1855.5682
$ws.Range("J34").Value = 4405.375
$ws.Range("K34").Value = 1855.5682
$ws.Range("L34").Value = 4405.375
$ws.Range("M34").Value = -1653.5682
$ws.Range("N34").Value = -4809.375
$ws.Range("H41").Value = 9799.6
$ws.Range("H50").Value = 9163.666999999999
$ws.Range("J50").Value = 9163.666999999999
$ws.Range("L50").Value = 9163.666999999999
$ws.Range("N50").Value = -10413.667
$ws.Range("H51").Value = 8585.637000000001
$ws.Range("J51").Value = 8872.75
$ws.Range("L51").Value = 8872.75
$ws.Range("N51").Value = -10344.75
$ws.Range("H59").Value = 15596.75
$ws.Range("J59").Value = 15596.75
$ws.Range("L59").Value = 15596.75
$ws.Range("N59").Value = -17886.75
$ws.Range("H60").Value = 6798
$ws.Range("J60").Value = 8064
$ws.Range("L60").Value = 8064
$ws.Range("N60").Value = -9086
$ws.Range("H61").Value = 8585.637000000001
$ws.Range("J61").Value = 8872.75
$ws.Range("L61").Value = 8872.75
$ws.Range("M61").Value = -7472
$ws.Range("N61").Value = -9568.75
$ws.Range("H62").Value = 6608.2
$ws.Range("I62").Value = 6863.8335
$ws.Range("J62").Value = 6224.75
$ws.Range("K62").Value = 6863.8335
$ws.Range("L62").Value = 6224.75
$ws.Range("M62").Value = -6239.8335
$ws.Range("N62").Value = -7472.75
$ws.Range("H65").Value = 6608.2
$ws.Range("I65").Value = 6863.8335
$ws.Range("J65").Value = 6224.75
$ws.Range("K65").Value = 34319.1675
$ws.Range("L65").Value = 31123.75
$ws.Range("M65").Value = -31199.1675
$ws.Range("N65").Value = -37363.75
$ws.Range("H74").Value = 13897
$ws.Range("J74").Value = 13897
$ws.Range("L74").Value = 13897
$ws.Range("N74").Value = -15645
$ws.Range("H77").Value = 13897
$ws.Range("J77").Value = 13897
$ws.Range("L77").Value = 41691
$ws.Range("N77").Value = -50427
$ws.Range("H141").Value = 66166.25
$ws.Range("J141").Value = 66166.25
$ws.Range("L141").Value = 66166.25
$ws.Range("N141").Value = -76526.25

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 116.5
$ws.Range("I14").Value = 116.5
$ws.Range("K14").Value = 349.5
$ws.Range("M14").Value = -176.5
$ws.Range("H23").Value = 118.05882
$ws.Range("I23").Value = 67.666664
$ws.Range("K23").Value = 202.999992
$ws.Range("M23").Value = 32.00000800000001
$ws.Range("H122").Value = 29144.139
$ws.Range("I122").Value = 578.2
$ws.Range("J122").Value = 33751.547
$ws.Range("K122").Value = 5203.8
$ws.Range("L122").Value = 303763.923
$ws.Range("M122").Value = -2753.8
$ws.Range("N122").Value = -308663.923
$ws.Range("H131").Value = 88690.25999999999
$ws.Range("J131").Value = 78869.69500000001
$ws.Range("L131").Value = 236609.085
$ws.Range("N131").Value = -246689.085
$ws.Range("H137").Value = 86545.836
$ws.Range("I137").Value = 3305
$ws.Range("K137").Value = 9915
$ws.Range("M137").Value = -4815

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4137.35
$ws.Range("I70").Value = 3999.6667
$ws.Range("J70").Value = 4196.357
$ws.Range("K70").Value = 3999.6667
$ws.Range("L70").Value = 4196.357
$ws.Range("M70").Value = -3729.6667
$ws.Range("N70").Value = -4736.357
$ws.Range("H73").Value = 4137.35
$ws.Range("I73").Value = 3999.6667
$ws.Range("J73").Value = 4196.357
$ws.Range("K73").Value = 3999.6667
$ws.Range("L73").Value = 4196.357
$ws.Range("M73").Value = -3063.6667
$ws.Range("N73").Value = -6068.357

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2324.8103
$ws.Range("I132").Value = 1521.2162
$ws.Range("J132").Value = 3740.6667
$ws.Range("K132").Value = 4563.6486
$ws.Range("L132").Value = 11222.0001
$ws.Range("M132").Value = -2033.6486
$ws.Range("N132").Value = -16282.0001
$ws.Range("H136").Value = 4683.162
$ws.Range("I136").Value = 2618.348
$ws.Range("J136").Value = 8075.357
$ws.Range("K136").Value = 7855.044
$ws.Range("L136").Value = 24226.071
$ws.Range("M136").Value = -5305.044
$ws.Range("N136").Value = -29326.071

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 864.5
$ws.Range("I126").Value = 838.6923
$ws.Range("K126").Value = 2516.0769
$ws.Range("M126").Value = -46.07690000000002
